$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.017899151126508
$ws.Cells.Item(2, 4).Value = 1.020544088705634
$ws.Cells.Item(2, 5).Value = 1.019234875372977
$ws.Cells.Item(2, 6).Value = 1.016269099976557
$ws.Cells.Item(2, 9).Value = 1.027402073278272
$ws.Cells.Item(2, 10).Value = 1.023110646044167
$ws.Cells.Item(2, 11).Value = 1.023384060977011
$ws.Cells.Item(2, 12).Value = 1.022078725679131
$ws.Cells.Item(2, 13).Value = 1.01912177548163
$ws.Cells.Item(2, 14).Value = 1.011782814275983

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.018843382827272
$ws.Cells.Item(3, 4).Value = 1.021400219807826
$ws.Cells.Item(3, 5).Value = 1.020035439859249
$ws.Cells.Item(3, 6).Value = 1.017854826511735
$ws.Cells.Item(3, 9).Value = 1.027493721362179
$ws.Cells.Item(3, 10).Value = 1.023690992705011
$ws.Cells.Item(3, 11).Value = 1.024046317098267
$ws.Cells.Item(3, 12).Value = 1.022685286405229
$ws.Cells.Item(3, 13).Value = 1.020510686138882
$ws.Cells.Item(3, 14).Value = 1.011976949331766

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.019454310653812
$ws.Cells.Item(4, 4).Value = 1.02195443230348
$ws.Cells.Item(4, 5).Value = 1.020553803289768
$ws.Cells.Item(4, 6).Value = 1.018880741203984
$ws.Cells.Item(4, 9).Value = 1.027551424534012
$ws.Cells.Item(4, 10).Value = 1.02406588396398
$ws.Cells.Item(4, 11).Value = 1.024474443555548
$ws.Cells.Item(4, 12).Value = 1.023077467842818
$ws.Cells.Item(4, 13).Value = 1.021408783869636
$ws.Cells.Item(4, 14).Value = 1.01210229408886

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.019711132102337
$ws.Cells.Item(5, 4).Value = 1.022187480023585
$ws.Cells.Item(5, 5).Value = 1.0207718053412
$ws.Cells.Item(5, 6).Value = 1.019312003953308
$ws.Cells.Item(5, 9).Value = 1.027575299816638
$ws.Cells.Item(5, 10).Value = 1.024223336969571
$ws.Cells.Item(5, 11).Value = 1.024654332642535
$ws.Cells.Item(5, 12).Value = 1.023242267976592
$ws.Cells.Item(5, 13).Value = 1.021786200453866
$ws.Cells.Item(5, 14).Value = 1.012154923455355

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.019754252817948
$ws.Cells.Item(6, 4).Value = 1.02222661308265
$ws.Cells.Item(6, 5).Value = 1.020808413644014
$ws.Cells.Item(6, 6).Value = 1.019384413166749
$ws.Cells.Item(6, 9).Value = 1.027579286106982
$ws.Cells.Item(6, 10).Value = 1.024249765150348
$ws.Cells.Item(6, 11).Value = 1.02468453122164
$ws.Cells.Item(6, 12).Value = 1.02326993436133
$ws.Cells.Item(6, 13).Value = 1.021849562071874
$ws.Cells.Item(6, 14).Value = 1.012163756310542

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.019457742363969
$ws.Cells.Item(7, 4).Value = 1.021957546075773
$ws.Cells.Item(7, 5).Value = 1.020556715921972
$ws.Cells.Item(7, 6).Value = 1.018886503875825
$ws.Cells.Item(7, 9).Value = 1.027551745062483
$ws.Cells.Item(7, 10).Value = 1.024067988453211
$ws.Cells.Item(7, 11).Value = 1.02447684761685
$ws.Cells.Item(7, 12).Value = 1.023079670197135
$ws.Cells.Item(7, 13).Value = 1.021413827487194
$ws.Cells.Item(7, 14).Value = 1.012102997582765

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.018218269751783
$ws.Cells.Item(8, 4).Value = 1.020833372679058
$ws.Cells.Item(8, 5).Value = 1.019505358253153
$ws.Cells.Item(8, 6).Value = 1.01680503947795
$ws.Cells.Item(8, 9).Value = 1.027433377158695
$ws.Cells.Item(8, 10).Value = 1.02330690773343
$ws.Cells.Item(8, 11).Value = 1.023607955785741
$ws.Cells.Item(8, 12).Value = 1.022283778533758
$ws.Cells.Item(8, 13).Value = 1.019591295683441
$ws.Cells.Item(8, 14).Value = 1.011848479738909

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.016033746361295
$ws.Cells.Item(9, 4).Value = 1.018854273722213
$ws.Cells.Item(9, 5).Value = 1.01765539074711
$ws.Cells.Item(9, 6).Value = 1.013135781918462
$ws.Cells.Item(9, 9).Value = 1.027212560688886
$ws.Cells.Item(9, 10).Value = 1.021960947183194
$ws.Cells.Item(9, 11).Value = 1.022073816602501
$ws.Cells.Item(9, 12).Value = 1.020878995687234
$ws.Cells.Item(9, 13).Value = 1.016374795255028
$ws.Cells.Item(9, 14).Value = 1.011397894034774

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.014577106309012
$ws.Cells.Item(10, 4).Value = 1.017536123836003
$ws.Cells.Item(10, 5).Value = 1.016423890336624
$ws.Cells.Item(10, 6).Value = 1.010688278942014
$ws.Cells.Item(10, 9).Value = 1.027057136662496
$ws.Cells.Item(10, 10).Value = 1.021060384087702
$ws.Cells.Item(10, 11).Value = 1.021049015286097
$ws.Cells.Item(10, 12).Value = 1.019940921769979
$ws.Cells.Item(10, 13).Value = 1.014226819818995
$ws.Cells.Item(10, 14).Value = 1.011096099374314

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.013946292211651
$ws.Cells.Item(11, 4).Value = 1.016965648658746
$ws.Cells.Item(11, 5).Value = 1.015891070837486
$ws.Cells.Item(11, 6).Value = 1.009628087928411
$ws.Cells.Item(11, 9).Value = 1.026987891084258
$ws.Cells.Item(11, 10).Value = 1.020669656655092
$ws.Cells.Item(11, 11).Value = 1.020604780376395
$ws.Cells.Item(11, 12).Value = 1.01953435819036
$ws.Cells.Item(11, 13).Value = 1.013295791202978
$ws.Cells.Item(11, 14).Value = 1.010965086224629

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.013711967405235
$ws.Cells.Item(12, 4).Value = 1.016753792765448
$ws.Cells.Item(12, 5).Value = 1.01569322255211
$ws.Cells.Item(12, 6).Value = 1.009234217867931
$ws.Cells.Item(12, 9).Value = 1.026961878009907
$ws.Cells.Item(12, 10).Value = 1.020524406041568
$ws.Cells.Item(12, 11).Value = 1.020439698130518
$ws.Cells.Item(12, 12).Value = 1.019383286563693
$ws.Cells.Item(12, 13).Value = 1.012949818428124
$ws.Cells.Item(12, 14).Value = 1.010916371909964

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.013762231438314
$ws.Cells.Item(13, 4).Value = 1.016799234582388
$ws.Cells.Item(13, 5).Value = 1.015735658753672
$ws.Cells.Item(13, 6).Value = 1.009318707517332
$ws.Cells.Item(13, 9).Value = 1.02696747112182
$ws.Cells.Item(13, 10).Value = 1.020555568091376
$ws.Cells.Item(13, 11).Value = 1.020475112166272
$ws.Cells.Item(13, 12).Value = 1.019415694466316
$ws.Cells.Item(13, 13).Value = 1.013024037523188
$ws.Cells.Item(13, 14).Value = 1.010926823567514

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.013926923104856
$ws.Cells.Item(14, 4).Value = 1.016948135686891
$ws.Cells.Item(14, 5).Value = 1.015874715317241
$ws.Cells.Item(14, 6).Value = 1.009595531931799
$ws.Cells.Item(14, 9).Value = 1.026985746794877
$ws.Cells.Item(14, 10).Value = 1.020657652583059
$ws.Cells.Item(14, 11).Value = 1.020591136129061
$ws.Cells.Item(14, 12).Value = 1.019521871701211
$ws.Cells.Item(14, 13).Value = 1.013267196024469
$ws.Cells.Item(14, 14).Value = 1.010961060508248

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.01402839342445
$ws.Cells.Item(15, 4).Value = 1.01703988440205
$ws.Cells.Item(15, 5).Value = 1.01596040122308
$ws.Cells.Item(15, 6).Value = 1.009766083383687
$ws.Cells.Item(15, 9).Value = 1.026996968328223
$ws.Cells.Item(15, 10).Value = 1.020720534691153
$ws.Cells.Item(15, 11).Value = 1.020662612562385
$ws.Cells.Item(15, 12).Value = 1.019587283603858
$ws.Cells.Item(15, 13).Value = 1.013416994307481
$ws.Cells.Item(15, 14).Value = 1.010982148365473

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.014618969653766
$ws.Cells.Item(16, 4).Value = 1.017573990560583
$ws.Cells.Item(16, 5).Value = 1.016459260857626
$ws.Cells.Item(16, 6).Value = 1.010758631183744
$ws.Cells.Item(16, 9).Value = 1.027061691279421
$ws.Cells.Item(16, 10).Value = 1.021086298978056
$ws.Cells.Item(16, 11).Value = 1.021078487376897
$ws.Cells.Item(16, 12).Value = 1.01996789623006
$ws.Cells.Item(16, 13).Value = 1.014288588790993
$ws.Cells.Item(16, 14).Value = 1.01110478724655

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.014989400966728
$ws.Cells.Item(17, 4).Value = 1.017909099574358
$ws.Cells.Item(17, 5).Value = 1.016772297177296
$ws.Cells.Item(17, 6).Value = 1.011381116615256
$ws.Cells.Item(17, 9).Value = 1.027101769414257
$ws.Cells.Item(17, 10).Value = 1.021315524900634
$ws.Cells.Item(17, 11).Value = 1.021339223606862
$ws.Cells.Item(17, 12).Value = 1.020206545037072
$ws.Cells.Item(17, 13).Value = 1.014835060870132
$ws.Cells.Item(17, 14).Value = 1.011181625932375

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.015205459659098
$ws.Cells.Item(18, 4).Value = 1.018104591203798
$ws.Cells.Item(18, 5).Value = 1.016954927275224
$ws.Cells.Item(18, 6).Value = 1.011744163120368
$ws.Cells.Item(18, 9).Value = 1.027124958544882
$ws.Cells.Item(18, 10).Value = 1.021449153483865
$ws.Cells.Item(18, 11).Value = 1.02149125944899
$ws.Cells.Item(18, 12).Value = 1.020345709016904
$ws.Cells.Item(18, 13).Value = 1.0151537183855
$ws.Cells.Item(18, 14).Value = 1.01122641240553

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.015279128809744
$ws.Cells.Item(19, 4).Value = 1.0181712535808
$ws.Cells.Item(19, 5).Value = 1.017017206417911
$ws.Cells.Item(19, 6).Value = 1.011867946239436
$ws.Cells.Item(19, 9).Value = 1.027132833585246
$ws.Cells.Item(19, 10).Value = 1.021494704659249
$ws.Cells.Item(19, 11).Value = 1.0215430917418
$ws.Cells.Item(19, 12).Value = 1.020393154260372
$ws.Cells.Item(19, 13).Value = 1.015262357208804
$ws.Cells.Item(19, 14).Value = 1.011241677968512

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.014949657984323
$ws.Cells.Item(20, 4).Value = 1.017873142630055
$ws.Cells.Item(20, 5).Value = 1.016738707068032
$ws.Cells.Item(20, 6).Value = 1.011314333873243
$ws.Cells.Item(20, 9).Value = 1.02709748883005
$ws.Cells.Item(20, 10).Value = 1.021290938889502
$ws.Cells.Item(20, 11).Value = 1.021311253958787
$ws.Cells.Item(20, 12).Value = 1.020180943984857
$ws.Cells.Item(20, 13).Value = 1.014776438980588
$ws.Cells.Item(20, 14).Value = 1.011173385202525

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.013878425858781
$ws.Cells.Item(21, 4).Value = 1.016904286808198
$ws.Cells.Item(21, 5).Value = 1.015833764839893
$ws.Cells.Item(21, 6).Value = 1.009514015971745
$ws.Cells.Item(21, 9).Value = 1.026980373128984
$ws.Cells.Item(21, 10).Value = 1.020627594477481
$ws.Cells.Item(21, 11).Value = 1.0205569719972
$ws.Cells.Item(21, 12).Value = 1.019490606692743
$ws.Cells.Item(21, 13).Value = 1.013195596019978
$ws.Cells.Item(21, 14).Value = 1.010950979968305

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.013204827947955
$ws.Cells.Item(22, 4).Value = 1.016295383220564
$ws.Cells.Item(22, 5).Value = 1.015265165620469
$ws.Cells.Item(22, 6).Value = 1.008381688287256
$ws.Cells.Item(22, 9).Value = 1.026905047440057
$ws.Cells.Item(22, 10).Value = 1.020209846204252
$ws.Cells.Item(22, 11).Value = 1.020082298974185
$ws.Cells.Item(22, 12).Value = 1.019056241320111
$ws.Cells.Item(22, 13).Value = 1.012200802172401
$ws.Cells.Item(22, 14).Value = 1.010810854537298

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.013561921896888
$ws.Cells.Item(23, 4).Value = 1.01661815039164
$ws.Cells.Item(23, 5).Value = 1.015566555303895
$ws.Cells.Item(23, 6).Value = 1.008981996502686
$ws.Cells.Item(23, 9).Value = 1.026945139186123
$ws.Cells.Item(23, 10).Value = 1.020431366742184
$ws.Cells.Item(23, 11).Value = 1.020333972550967
$ws.Cells.Item(23, 12).Value = 1.019286537238445
$ws.Cells.Item(23, 13).Value = 1.012728244250009
$ws.Cells.Item(23, 14).Value = 1.010885165212439

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.014967616150681
$ws.Cells.Item(24, 4).Value = 1.017889389938134
$ws.Cells.Item(24, 5).Value = 1.016753884864814
$ws.Cells.Item(24, 6).Value = 1.011344510237738
$ws.Cells.Item(24, 9).Value = 1.027099423621943
$ws.Cells.Item(24, 10).Value = 1.021302048481274
$ws.Cells.Item(24, 11).Value = 1.021323892384466
$ws.Cells.Item(24, 12).Value = 1.02019251210896
$ws.Cells.Item(24, 13).Value = 1.014802927967044
$ws.Cells.Item(24, 14).Value = 1.011177108933174

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.016598548387758
$ws.Cells.Item(25, 4).Value = 1.019365698980129
$ws.Cells.Item(25, 5).Value = 1.018133334130978
$ws.Cells.Item(25, 6).Value = 1.014084581459485
$ws.Cells.Item(25, 9).Value = 1.027271095695447
$ws.Cells.Item(25, 10).Value = 1.022309483938433
$ws.Cells.Item(25, 11).Value = 1.022470788406222
$ws.Cells.Item(25, 12).Value = 1.021242439827313
$ws.Cells.Item(25, 13).Value = 1.017206957359882
$ws.Cells.Item(25, 14).Value = 1.011514628936236

